$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.152.97"
$ws.Range("E2").Value = "  +1.37%  "
$ws.Range("D3").Value = "2.585.77"
$ws.Range("E3").Value = "  -0.58%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'521.38"
$ws.Range("E5").Value = "  +0.17%  "
$ws.Range("D6").Value = "'138.96"
$ws.Range("E6").Value = "  -3.89%  "
$ws.Range("E7").Value = "  +0.29%  "
$ws.Range("E8").Value = "  -0.63%  "
$ws.Range("D9").Value = "2.596.58"
$ws.Range("E9").Value = "  -1.13%  "
$ws.Range("E10").Value = "  -2.45%  "
$ws.Range("E11").Value = "  -1.06%  "
$ws.Range("D12").Value = "'0.329"
$ws.Range("E12").Value = "  +0.45%  "
$ws.Range("E13").Value = "  +3.15%  "
$ws.Range("D14").Value = "3.048.73"
$ws.Range("E14").Value = "  -0.38%  "
$ws.Range("D15").Value = "59.076.03"
$ws.Range("E15").Value = "  +1.34%  "
$ws.Range("D16").Value = "'20.55"
$ws.Range("E16").Value = "  -0.11%  "
$ws.Range("D17").Value = "2.603.31"
$ws.Range("E17").Value = "  +0.52%  "
$ws.Range("E18").Value = "  -1.56%  "
$ws.Range("D19").Value = "'338.27"
$ws.Range("E19").Value = "  -1.11%  "
$ws.Range("E20").Value = "  -1.06%  "
$ws.Range("D21").Value = "'10.06"
$ws.Range("E21").Value = "  -2.87%  "
$ws.Range("D22").Value = "'6.46"
$ws.Range("E22").Value = "  +0.83%  "
$ws.Range("E23").Value = "  +0.02%  "
$ws.Range("D24").Value = "'66.19"
$ws.Range("E24").Value = "  -0.35%  "
$ws.Range("E25").Value = "  +1.30%  "
$ws.Range("E26").Value = "  -0.39%  "
$ws.Range("E27").Value = "  +0.31%  "
$ws.Range("E28").Value = "  -0.65%  "
$ws.Range("D29").Value = "'0.999"
$ws.Range("E29").Value = "  +0.11%  "
$ws.Range("D30").Value = "0.0₃0723"
$ws.Range("E30").Value = "  -4.29%  "
$ws.Range("E31").Value = "  -6.38%  "
$ws.Range("E32").Value = "  -0.66%  "
$ws.Range("D33").Value = "'18.69"
$ws.Range("E33").Value = "  -0.94%  "
$ws.Range("D34").Value = "'149.50"
$ws.Range("E34").Value = "  -0.14%  "
$ws.Range("D35").Value = "'3.97"
$ws.Range("E35").Value = "  -2.27%  "
$ws.Range("E36").Value = "  -3.16%  "
$ws.Range("D37").Value = "'36.77"
$ws.Range("E37").Value = "  +1.64%  "
$ws.Range("D38").Value = "'1.46"
$ws.Range("E38").Value = "  -0.65%  "
$ws.Range("D39").Value = "'0.823"
$ws.Range("E39").Value = "  -2.98%  "
$ws.Range("D40").Value = "'0.815"
$ws.Range("E40").Value = "  -7.45%  "
$ws.Range("E41").Value = "  -1.54%  "
$ws.Range("E42").Value = "  +0.26%  "
$ws.Range("D43").Value = "'272.48"
$ws.Range("E43").Value = "  -1.71%  "
$ws.Range("D44").Value = "'10.76"
$ws.Range("E44").Value = "  +0.93%  "
$ws.Range("E45").Value = "  -1.36%  "
$ws.Range("D46").Value = "'0.0951"
$ws.Range("E46").Value = "  -0.63%  "
$ws.Range("E47").Value = "  -1.69%  "
$ws.Range("D48").Value = "'18.38"
$ws.Range("E48").Value = "  -3.19%  "
$ws.Range("D49").Value = "1.972.08"
$ws.Range("E49").Value = "  -0.73%  "
$ws.Range("E50").Value = "  -0.72%  "
$ws.Range("D51").Value = "'4.50"
$ws.Range("E51").Value = "  -5.01%  "
